$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.012.17"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").Value = "2.461.39"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'517.09"
$ws.Range("E5").Value = "  -3.59%  "
$ws.Range("D6").Value = "'131.02"
$ws.Range("E6").Value = "  -3.70%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'0.556"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "2.464.38"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("D10").Value = "'0.0985"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "'5.26"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "'0.340"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "2.899.17"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "57.950.95"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "'22.19"
$ws.Range("E16").Value = "  -3.26%  "
$ws.Range("D17").Value = "'0.0000135"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").Value = "2.460.44"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").Value = "'10.67"
$ws.Range("E19").Value = "  -3.71%  "
$ws.Range("D20").Value = "'319.01"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("E23").Value = "  -3.54%  "
$ws.Range("D24").Value = "'64.20"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").Value = "'0.406"
$ws.Range("E25").Value = "  -3.04%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("D28").Value = "'7.29"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").Value = "0.0₃0733"
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("D30").Value = "'165.91"
$ws.Range("E30").Value = "  -2.80%  "
$ws.Range("E31").Value = "  -4.14%  "
$ws.Range("D32").Value = "'6.21"
$ws.Range("E32").Value = "  -5.83%  "
$ws.Range("D33").Value = "'1.15"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "'18.01"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("E37").Value = "  -7.74%  "
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").Value = "'1.46"
$ws.Range("D40").Value = "'0.782"
$ws.Range("E40").Value = "  -3.20%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.42"
$ws.Range("E41").Value = "  -4.32%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'270.63"
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'4.94"
$ws.Range("E43").Value = "  -4.45%  "
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("D45").Value = "'124.43"
$ws.Range("E45").Value = "  -4.55%  "
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("E48").Value = "  -4.82%  "
$ws.Range("D49").Value = "'16.67"
$ws.Range("E49").Value = "  -3.74%  "
$ws.Range("D50").Value = "1.716.73"
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("E51").Value = "  -2.38%  "
